# Update "想去人数" (want-to-go count) figures on the 展览 (Exhibition),
# 演出 (Performance) and 全部类型 (All Types) sheets, matching the
# regenerated site data output.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet (column F) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 29
$ws1.Range("F4").Value = 43
$ws1.Range("F5").Value = 4951
$ws1.Range("F7").Value = 78
$ws1.Range("F8").Value = 277
$ws1.Range("F9").Value = 40

# --- 演出 sheet (column F) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 126

# --- 全部类型 sheet (column F) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 126
$ws4.Range("F7").Value = 29
$ws4.Range("F8").Value = 43
$ws4.Range("F9").Value = 4951
$ws4.Range("F11").Value = 78
$ws4.Range("F13").Value = 278
$ws4.Range("F14").Value = 40
